# Inserts a new "2022-10-11" Kiwi market-day entry (Especial/Primera/Segunda)
# right before the current row 283, shifting the remaining rows (old 283-392)
# down to 286-395, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 283..285 (existing rows shift down by 3)
$ws.Rows("283:285").Insert()

# Common (constant) column values used throughout this data block
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId = 100101
$producto  = "Berries"
$categoriaId = 100101007
$categoria = "Kiwi"
$variedad  = "Hayward"
$unidad    = "$/caja 15 kilos"
$origen    = "Región de O'Higgins"
$kgUnidad  = 15

$fecha = 44845  # 2022-10-11

# Row 283: Especial
$r = 283
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 17000
$ws.Cells.Item($r, 15).Value = 17000
$ws.Cells.Item($r, 16).Value = 17000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 1133
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 284: Primera
$r = 284
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 14000
$ws.Cells.Item($r, 15).Value = 14000
$ws.Cells.Item($r, 16).Value = 14000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 933
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 285: Segunda
$r = 285
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 12500
$ws.Cells.Item($r, 15).Value = 12500
$ws.Cells.Item($r, 16).Value = 12500
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 833
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Ensure the date column keeps the date number format used elsewhere (D column)
$ws.Range("D283:D285").NumberFormat = "YYYY-MM-DD HH:MM:SS"
